$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.842.60"
$ws.Cells.Item(2, 5).Value = "  +0.47%  "
$ws.Cells.Item(3, 4).Value = "3.496.64"
$ws.Cells.Item(3, 5).Value = "  +0.26%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'594.25"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.62%  "
$ws.Cells.Item(6, 4).Value = "'172.68"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +2.69%  "
$ws.Cells.Item(7, 5).Value = "  +0.00%  "
$ws.Cells.Item(8, 5).Value = "  -1.60%  "
$ws.Cells.Item(9, 5).Value = "  +4.23%  "
$ws.Cells.Item(10, 5).Value = "  -2.07%  "
$ws.Cells.Item(11, 4).Value = "'0.431"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.57%  "
$ws.Cells.Item(12, 4).Value = "4.098.49"
$ws.Cells.Item(12, 5).Value = "  +0.17%  "
$ws.Cells.Item(13, 5).Value = "  +0.13%  "
$ws.Cells.Item(14, 4).Value = "'29.27"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +4.37%  "
$ws.Cells.Item(15, 4).Value = "66.848.24"
$ws.Cells.Item(15, 5).Value = "  +0.44%  "
$ws.Cells.Item(16, 5).Value = "  +0.75%  "
$ws.Cells.Item(17, 4).Value = "3.500.63"
$ws.Cells.Item(17, 5).Value = "  -0.14%  "
$ws.Cells.Item(18, 4).Value = "'6.26"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.10%  "
$ws.Cells.Item(19, 4).Value = "'14.21"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +2.08%  "
$ws.Cells.Item(20, 4).Value = "'392.95"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.08%  "
$ws.Cells.Item(21, 5).Value = "  +0.24%  "
$ws.Cells.Item(22, 4).Value = "'73.26"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.56%  "
$ws.Cells.Item(23, 5).Value = "  +0.15%  "
$ws.Cells.Item(24, 4).Value = "'0.534"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.46%  "
$ws.Cells.Item(25, 5).Value = "  +0.05%  "
$ws.Cells.Item(26, 4).Value = "'10.23"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.43%  "
$ws.Cells.Item(27, 5).Value = "  +0.58%  "
$ws.Cells.Item(28, 4).Value = "'0.998"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.29%  "
$ws.Cells.Item(29, 4).Value = "'6.12"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -2.79%  "
$ws.Cells.Item(30, 5).Value = "  -1.76%  "
$ws.Cells.Item(31, 5).Value = "  +0.03%  "
$ws.Cells.Item(32, 4).Value = "'23.63"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.42%  "
$ws.Cells.Item(33, 5).Value = "  +0.10%  "
$ws.Cells.Item(34, 5).Value = "  +0.99%  "
$ws.Cells.Item(35, 4).Value = "'163.28"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.35%  "
$ws.Cells.Item(36, 5).Value = "  -1.37%  "
$ws.Cells.Item(37, 5).Value = "  -0.74%  "
$ws.Cells.Item(38, 4).Value = "'6.85"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.40%  "
$ws.Cells.Item(39, 5).Value = "  +0.30%  "
$ws.Cells.Item(40, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(40, 4).Value = "'27.22"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.58%  "
$ws.Cells.Item(41, 2).Value = "Maker"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(41, 4).Value = "2.842.00"
$ws.Cells.Item(41, 5).Value = "  +2.25%  "
$ws.Cells.Item(42, 5).Value = "  -0.64%  "
$ws.Cells.Item(43, 4).Value = "'26.01"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.95%  "
$ws.Cells.Item(44, 4).Value = "'42.63"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.94%  "
$ws.Cells.Item(45, 4).Value = "'2.53"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.09%  "
$ws.Cells.Item(46, 4).Value = "'0.0302"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.53%  "
$ws.Cells.Item(47, 4).Value = "'337.75"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.19%  "
$ws.Cells.Item(48, 4).Value = "'34.62"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.78%  "
$ws.Cells.Item(49, 5).Value = "  -0.96%  "
$ws.Cells.Item(50, 4).Value = "'6.43"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.88%  "
$ws.Cells.Item(51, 5).Value = "  -2.81%  "
